$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "HPO Terms"

$ws.Range("H2").Select()
